$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 577; existing rows 577:693 shift down to 578:694
$ws.Rows.Item(577).Insert()

# Populate the newly inserted row 577 with the new record's data
$ws.Cells.Item(577, 1).Value = 9
$ws.Cells.Item(577, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(577, 3).Value = "Metropolitana"
$ws.Cells.Item(577, 4).Value = 45258
$ws.Cells.Item(577, 5).Value = 13
$ws.Cells.Item(577, 6).Value = 100112012
$ws.Cells.Item(577, 7).Value = "Espinaca"
$ws.Cells.Item(577, 8).Value = "Sin especificar"
$ws.Cells.Item(577, 9).Value = "Primera"
$ws.Cells.Item(577, 10).Value = 130
$ws.Cells.Item(577, 11).Value = 15000
$ws.Cells.Item(577, 12).Value = 17000
$ws.Cells.Item(577, 13).Value = 16231
$ws.Cells.Item(577, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(577, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(577, 16).Value = 1623
$ws.Cells.Item(577, 17).Value = 10
$ws.Cells.Item(577, 18).Value = "Hortaliza"
